$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Verified"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Interior.Color = 65535

$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 1
